# film_info.xlsx edit script
# - Quality (col F) screenshots bumped from 720 -> 1080 for most rows
#   (rows 15 and 23 intentionally stay at 720)
# - Pulp Fiction (row 19) clip timestamps nudged +2s each, trimming the
#   counted duration by 4s (the first 4s of black silence were removed)
# - A threaded-comment reply was added on D19 explaining the manual edit
# - Active selection left on I26

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reply to the existing threaded comment on D19 -------------------------
$ct = $ws.Range("D19").CommentThreaded()
$ct.AddReply("I manually removed the first four seconds of the original ratings file.")

# --- Pulp Fiction (row 19): drop first 4 seconds of ratings ----------------
$ws.Range("C19").Value = "00:14:25"
$ws.Range("D19").Value = "00:18:18"
$ws.Range("E19").Value = 229

# --- Quality column (F): 720 -> 1080, skipping rows 15 and 23 --------------
$qualityRows = @(4,5,6,7,8,9,11,14,16,17,18,19,20,22)
foreach ($r in $qualityRows) {
    $ws.Cells.Item($r, 6).Value = 1080
}

# --- Leave the selection where the author last left it ---------------------
$ws.Range("I26").Select()
